$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column E; everything from E onward shifts right by one.
$ws.Columns("E:E").Insert()

# The insert pushed the old empty trailing column (AA) out to AB; drop it so the
# used range stays at A1:AA2.
$ws.Columns("AB:AB").Delete()

# Populate the newly inserted column E.
$ws.Range("E1").Value = "DropDownOption"
$ws.Range("E2").Value = "Incident Report"

# Fix up labels that moved right by one column and were also renamed.
$ws.Range("P1").Value = "SendOneTimeEmailNotificationOnSave"
$ws.Range("S1").Value = "NotifyCCRecipientsOnCreateEditandClose"
$ws.Range("U1").Value = "Creator_OnIRCreateandEdit"
$ws.Range("V1").Value = "Creator_OnIRClose"
$ws.Range("W1").Value = "Creator_OnAssigneeReassign"
$ws.Range("Y1").Value = "Assignee_OnIRCreateandEdit"
$ws.Range("Z1").Value = "Assignee_OnIRClose"
$ws.Range("AA1").Value = "Assignee_OnAssigneeReassign"

# Rename the worksheet.
$ws.Name = "testIR_GeneralInformation"

# Update the active selection.
$ws.Range("B1").Select()
